# Add a new row (row 12) of medical record data to Sheet1, matching the
# next entry appended to the Medical_Records table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(12, 1).Value = "Mahatab Ali"

# Column B ("Age") holds numeric-looking text ("10") in this sheet, just
# like the existing rows (e.g. B11 = "27"). Force it to be stored as text
# rather than being auto-converted to a number by Excel.
$bCell = $ws.Cells.Item(12, 2)
$bCell.NumberFormat = "@"
$bCell.Value = "10"
$bCell.Style = "Normal"

$ws.Cells.Item(12, 3).Value = "Male"
$ws.Cells.Item(12, 4).Value = "16-01-2025 23:05:10"
$ws.Cells.Item(12, 5).Value = "pain right groin for last 3 days with fever."
$ws.Cells.Item(12, 6).Value = "Anterior hip point(right) tender"
$ws.Cells.Item(12, 7).Value = "TC DC ESR, Hb, CRP"
$ws.Cells.Item(12, 8).Value = "Tablet CETIL 250 mg 1 tab BDPC X 5 days."
